# "Copy in EU-2024-develop branch"
# Update the "share of costs that must be covered" values on the
# SoCtMbCtbDP sheet from 0.95 to 1 for every technology row (B2:B25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCtMbCtbDP")

$ws.Range("B2:B25").Value = 1
